# Project Vision and Scope Document - "Customer or Market Needs" minor edit
#
# 1) The stray "_GoBack" bookmark sitting between "...it is located " and
#    "beside" is removed (Word drops/relocates _GoBack as the document is
#    edited). Removing it also renumbers the bookmarks that follow it
#    (_Toc18551420, _Toc18551421, ...) down by one.
# 2) A fresh "_GoBack" bookmark is dropped where the author's cursor last
#    was: right after "...check availability through" in the "Customer or
#    Market Needs" paragraph.
# 3) The word "the" is inserted so the sentence reads "...so that it will
#    enhance the user experience..." instead of "...enhance user
#    experience...".

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -----------------------------------
$oldGoBack = $d.Bookmarks("_GoBack")
$oldGoBack.Delete()

# --- 2. Re-create _GoBack at its new location ------------------------------
$rng = $d.Content
$rng.Find.Execute("check availability through", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPoint = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

# --- 3. Insert the missing word "the" before "user experience" ------------
$rng2 = $d.Content
$rng2.Find.Execute("user experience", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$theePoint = $d.Range($rng2.Start, $rng2.Start)
$theePoint.InsertBefore("the ")
